# Change the table style (tableStyleId) of the table on slide 16 from the
# custom "Table_0" style ({EE90A90D-B5BF-4095-B2B0-BE8BBF14C4B1}) to the
# built-in PowerPoint table style {2431BC8A-34E7-4E55-958C-8CDA5FE98859}.
#
# PowerPoint's Table object exposes the current style id through the
# read/write-looking `Style` property, but it cannot be assigned directly
# (doing so raises "Table styles cannot be assigned through a property -
# call Table.ApplyStyle(...) instead"). The supported way to change it is
# Table.ApplyStyle("{GUID}").

$p = $ppt.ActivePresentation

$targetStyleId = "{2431BC8A-34E7-4E55-958C-8CDA5FE98859}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($targetStyleId)
        }
    }
}
